$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "2026-02-08 23:18:25"
$ws.Cells.Item(3, 5).Value = "2026-02-08 23:18:27"
$ws.Cells.Item(4, 5).Value = "2026-02-08 23:18:30"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = "72%"
$ws.Cells.Item(4, 10).Value = "1002.9 hPa"
$ws.Cells.Item(4, 15).Value = "10.0 °C"
$ws.Cells.Item(5, 5).Value = "2026-02-08 23:18:32"
$ws.Cells.Item(6, 5).Value = "2026-02-08 23:18:35"
$ws.Cells.Item(6, 10).Value = "1002.8 hPa"
$ws.Cells.Item(6, 15).Value = "9.9 °C"
$ws.Cells.Item(7, 5).Value = "2026-02-08 23:18:37"
$ws.Cells.Item(7, 10).Value = "1003.1 hPa"
$ws.Cells.Item(8, 5).Value = "2026-02-08 23:18:40"
$ws.Cells.Item(8, 10).Value = "1003.1 hPa"
$ws.Cells.Item(9, 5).Value = "2026-02-08 23:18:42"
$ws.Cells.Item(10, 5).Value = "2026-02-08 23:18:45"
$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = "81%"
$ws.Cells.Item(11, 5).Value = "2026-02-08 23:18:47"
$ws.Cells.Item(12, 5).Value = "2026-02-08 23:18:50"
$ws.Cells.Item(13, 5).Value = "2026-02-08 23:18:52"
$ws.Cells.Item(13, 10).Value = "1004.4 hPa"
$ws.Cells.Item(14, 5).Value = "2026-02-08 23:18:54"
$ws.Cells.Item(14, 8).NumberFormat = "@"
$ws.Cells.Item(14, 8).Value = "77%"
$ws.Cells.Item(14, 15).Value = "11.0 °C"
$ws.Cells.Item(15, 5).Value = "2026-02-08 23:18:57"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "73%"
$ws.Cells.Item(16, 5).Value = "2026-02-08 23:18:59"
$ws.Cells.Item(17, 5).Value = "2026-02-08 23:19:02"
$ws.Cells.Item(18, 5).Value = "2026-02-08 23:19:04"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value = "72%"
$ws.Cells.Item(18, 10).Value = "1003.2 hPa"
$ws.Cells.Item(18, 15).Value = "10.1 °C"
$ws.Cells.Item(19, 5).Value = "2026-02-08 23:19:06"
$ws.Cells.Item(20, 5).Value = "2026-02-08 23:19:09"
$ws.Cells.Item(20, 9).Value = "10.2 mm"
$ws.Cells.Item(21, 5).Value = "2026-02-08 23:19:11"
$ws.Cells.Item(21, 10).Value = "1003.8 hPa"
$ws.Cells.Item(22, 5).Value = "2026-02-08 23:19:14"
$ws.Cells.Item(23, 5).Value = "2026-02-08 23:19:16"
$ws.Cells.Item(24, 5).Value = "2026-02-08 23:19:19"
$ws.Cells.Item(24, 10).Value = "1004.3 hPa"
$ws.Cells.Item(24, 15).Value = "8.3 °C"
$ws.Cells.Item(25, 5).Value = "2026-02-08 23:19:21"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = "81%"
$ws.Cells.Item(25, 15).Value = "-3.1 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-08 23:19:24"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "72%"
$ws.Cells.Item(26, 10).Value = "1002.3 hPa"
$ws.Cells.Item(26, 15).Value = "3.4 °C"
$ws.Cells.Item(27, 5).Value = "2026-02-08 23:19:26"
$ws.Cells.Item(28, 5).Value = "2026-02-08 23:19:29"
$ws.Cells.Item(28, 10).Value = "1002.8 hPa"
$ws.Cells.Item(28, 15).Value = "8.4 °C"
$ws.Cells.Item(29, 5).Value = "2026-02-08 23:19:31"
$ws.Cells.Item(30, 5).Value = "2026-02-08 23:19:34"
$ws.Cells.Item(30, 10).Value = "1003.1 hPa"
$ws.Cells.Item(31, 5).Value = "2026-02-08 23:19:36"
$ws.Cells.Item(31, 8).NumberFormat = "@"
$ws.Cells.Item(31, 8).Value = "78%"
$ws.Cells.Item(31, 10).Value = "1002.3 hPa"
$ws.Cells.Item(32, 5).Value = "2026-02-08 23:19:39"
$ws.Cells.Item(32, 8).NumberFormat = "@"
$ws.Cells.Item(32, 8).Value = "90%"
$ws.Cells.Item(33, 5).Value = "2026-02-08 23:19:41"
$ws.Cells.Item(33, 10).Value = "1003.9 hPa"
$ws.Cells.Item(34, 5).Value = "2026-02-08 23:19:44"
$ws.Cells.Item(35, 5).Value = "2026-02-08 23:19:46"
$ws.Cells.Item(35, 10).Value = "1005.3 hPa"
$ws.Cells.Item(36, 5).Value = "2026-02-08 23:19:49"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "75%"
$ws.Cells.Item(36, 10).Value = "1003.2 hPa"
$ws.Cells.Item(37, 5).Value = "2026-02-08 23:19:51"
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value = "80%"
$ws.Cells.Item(37, 10).Value = "1004.0 hPa"
$ws.Cells.Item(38, 5).Value = "2026-02-08 23:19:54"
$ws.Cells.Item(38, 8).NumberFormat = "@"
$ws.Cells.Item(38, 8).Value = "79%"
$ws.Cells.Item(39, 5).Value = "2026-02-08 23:19:56"
$ws.Cells.Item(39, 8).NumberFormat = "@"
$ws.Cells.Item(39, 8).Value = "87%"
$ws.Cells.Item(40, 5).Value = "2026-02-08 23:19:59"
$ws.Cells.Item(40, 8).NumberFormat = "@"
$ws.Cells.Item(40, 8).Value = "85%"
$ws.Cells.Item(40, 10).Value = "1004.5 hPa"
$ws.Cells.Item(40, 15).Value = "5.4 °C"
$ws.Cells.Item(41, 5).Value = "2026-02-08 23:20:01"
$ws.Cells.Item(41, 8).NumberFormat = "@"
$ws.Cells.Item(41, 8).Value = "66%"
$ws.Cells.Item(41, 10).Value = "1003.1 hPa"
$ws.Cells.Item(42, 5).Value = "2026-02-08 23:20:04"
$ws.Cells.Item(42, 8).NumberFormat = "@"
$ws.Cells.Item(42, 8).Value = "81%"
$ws.Cells.Item(43, 5).Value = "2026-02-08 23:20:06"
$ws.Cells.Item(44, 5).Value = "2026-02-08 23:20:08"
$ws.Cells.Item(45, 5).Value = "2026-02-08 23:20:11"
$ws.Cells.Item(45, 10).Value = "1005.3 hPa"
$ws.Cells.Item(46, 5).Value = "2026-02-08 23:20:13"
$ws.Cells.Item(46, 10).Value = "1004.9 hPa"
$ws.Cells.Item(46, 15).Value = "9.6 °C"
